# feat: Add zombie enemy
#
# Adds a new "Zombie" enemy to the game config workbook:
#   1. enemy-chances  - new spawn-chance row (level, tag, weight)
#   2. stuff-descriptor - new descriptor row (tag, icon, name, description,
#      color, exp, hp, melee_power, melee_skill, defense)
# and restores the selection/active-tab state recorded after the edit.

$wb = $excel.ActiveWorkbook

$wsEnemy = $wb.Worksheets.Item("enemy-chances")
$wsItem  = $wb.Worksheets.Item("item-chances")
$wsStuff = $wb.Worksheets.Item("stuff-descriptor")

# --- 1. enemy-chances: bump every existing enemy's spawn level by one ----
# (the new Zombie takes over the level-2 slot alongside Troll) and append
# the Zombie spawn-chance row (row 10).
$wsEnemy.Range("A2").Value2 = 1   # Orc: 0 -> 1
$wsEnemy.Range("A3").Value2 = 1   # Goblin: 0 -> 1
# A4 (Troll, level 2) is unchanged
$wsEnemy.Range("A5").Value2 = 4   # Gargoyle: 3 -> 4
$wsEnemy.Range("A6").Value2 = 5   # Troll: 4 -> 5
$wsEnemy.Range("A7").Value2 = 6   # Warlord: 5 -> 6
$wsEnemy.Range("A8").Value2 = 7   # Minotaur: 6 -> 7
$wsEnemy.Range("A9").Value2 = 7   # Warlord: 6 -> 7

$wsEnemy.Range("A10").Value2 = 2
$wsEnemy.Range("B10").Value2 = "Zombie"
$wsEnemy.Range("C10").Value2 = 30

# --- 2. stuff-descriptor: append the Zombie descriptor row (row 24) ------
$wsStuff.Range("A24").Value2 = "Zombie"
$wsStuff.Range("B24").Value2 = "zombie"
$wsStuff.Range("C24").Value2 = "Zombie"
$wsStuff.Range("D24").Value2 = "Shambling corpse. Once belonged to an adventurer like  you"
$wsStuff.Range("E24").Value2 = "#f0ddd7"
$wsStuff.Range("F24").Value2 = 50
$wsStuff.Range("G24").Value2 = 5
$wsStuff.Range("H24").Value2 = 3
$wsStuff.Range("I24").Value2 = 3
$wsStuff.Range("J24").Value2 = 1

# --- 3. Restore per-sheet selections, then leave enemy-chances active ----
# (selecting on a sheet makes it the active/tab-selected one in this
# runtime, so touch the non-active sheets first and finish on the sheet
# that should end up active.)
$wsStuff.Range("J24").Select()
$wsItem.Range("C12").Select()
$wsEnemy.Range("A12").Select()
